# LMS-2340 Finished implementing changes to templates and handlers.
# Apply the changes described by the diff:
#  - openbis-metadata!B2 gets the new value "/TEST/TEST/TEST"
#  - openbis-metadata sheet selection moves from B9 to C9
#  - Both sheets get fullCalcOnLoad enabled
#  - openbis-metadata sheet loses its explicit pageSetup

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("openbis-metadata")
$wsData = $wb.Worksheets.Item("openbis-data")

# Set the new shared string value on B2 of the metadata sheet.
$wsMeta.Range("B2").Value = "/TEST/TEST/TEST"

# Update the active cell / selection on the metadata sheet from B9 to C9.
$wsMeta.Range("C9").Select()

# Force both sheets to recalculate fully when the workbook is loaded.
$wsMeta.EnableCalculation = $true
$wsData.EnableCalculation = $true
$wb.Application.CalculateFullRebuild()

# Remove the explicit page setup on the metadata sheet.
$wsMeta.PageSetup.PaperSize = $null
